$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D (shifts the old D..J block -> E..K); this also
# updates the dimension, the merged-cell range and every existing cell's
# column letter automatically.
$ws.Columns("D").Insert()

# Give the new "password" column a Text number format *before* writing the
# values so the leading zeros in "00000000" are preserved as text instead
# of being coerced to the number 0.
$ws.Range("D1:D4").NumberFormat = "@"

$ws.Range("D1").Value = "密码"
$ws.Range("D2").Value = "00000000"
$ws.Range("D3").Value = "00000000"
$ws.Range("D4").Value = "00000000"

# Row 2 picks up an explicit (custom) row height matching the default.
$ws.Rows(2).RowHeight = 14.25

# Match the final UI selection left behind by the save.
$ws.Range("F17").Select() | Out-Null
